$d = $word.ActiveDocument

# N05 -> N04 and merge the trailing "s" run into the main text run so the
# paragraph ends up with a single run reading
# "N04: Controle de pagamentos de clientes".
$d.Content.Find.Execute(
    "N05: Controle de pagamentos de clientes",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "N04: Controle de pagamentos de clientes", 2)

# The section's page size now records an explicit portrait orientation.
$d.PageSetup.Orientation = 0
